$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-6
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = 4
